# Update rows 2-16 of the active sheet: the underlying daily price records
# were re-shuffled (weekly consolidation), so each row's Fecha (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M) and
# Precio $/Kg (P) values are rewritten to reflect the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44400
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 25000
$ws.Cells.Item(2, 12).Value = 25000
$ws.Cells.Item(2, 13).Value = 25000
$ws.Cells.Item(2, 16).Value = 1667

$ws.Cells.Item(3, 4).Value = 44754
$ws.Cells.Item(3, 10).Value = 90
$ws.Cells.Item(3, 11).Value = 25000
$ws.Cells.Item(3, 12).Value = 25000
$ws.Cells.Item(3, 13).Value = 25000
$ws.Cells.Item(3, 16).Value = 1667

$ws.Cells.Item(4, 4).Value = 44750
$ws.Cells.Item(4, 10).Value = 90
$ws.Cells.Item(4, 11).Value = 25000
$ws.Cells.Item(4, 12).Value = 25000
$ws.Cells.Item(4, 13).Value = 25000
$ws.Cells.Item(4, 16).Value = 1667

$ws.Cells.Item(5, 4).Value = 44740
$ws.Cells.Item(5, 10).Value = 90
$ws.Cells.Item(5, 11).Value = 25000
$ws.Cells.Item(5, 12).Value = 25000
$ws.Cells.Item(5, 13).Value = 25000
$ws.Cells.Item(5, 16).Value = 1667

$ws.Cells.Item(6, 4).Value = 44407
$ws.Cells.Item(6, 10).Value = 90
$ws.Cells.Item(6, 11).Value = 25000
$ws.Cells.Item(6, 12).Value = 25000
$ws.Cells.Item(6, 13).Value = 25000
$ws.Cells.Item(6, 16).Value = 1667

$ws.Cells.Item(7, 4).Value = 44757
$ws.Cells.Item(7, 10).Value = 80
$ws.Cells.Item(7, 11).Value = 25000
$ws.Cells.Item(7, 12).Value = 25000
$ws.Cells.Item(7, 13).Value = 25000
$ws.Cells.Item(7, 16).Value = 1667

$ws.Cells.Item(8, 4).Value = 44418
$ws.Cells.Item(8, 10).Value = 90
$ws.Cells.Item(8, 11).Value = 25000
$ws.Cells.Item(8, 12).Value = 25000
$ws.Cells.Item(8, 13).Value = 25000
$ws.Cells.Item(8, 16).Value = 1667

$ws.Cells.Item(9, 4).Value = 44775
$ws.Cells.Item(9, 10).Value = 120
$ws.Cells.Item(9, 11).Value = 24000
$ws.Cells.Item(9, 12).Value = 24000
$ws.Cells.Item(9, 13).Value = 24000
$ws.Cells.Item(9, 16).Value = 1600

$ws.Cells.Item(10, 4).Value = 44764
$ws.Cells.Item(10, 10).Value = 90
$ws.Cells.Item(10, 11).Value = 24000
$ws.Cells.Item(10, 12).Value = 24000
$ws.Cells.Item(10, 13).Value = 24000
$ws.Cells.Item(10, 16).Value = 1600

$ws.Cells.Item(11, 4).Value = 44778
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 24000
$ws.Cells.Item(11, 12).Value = 24000
$ws.Cells.Item(11, 13).Value = 24000
$ws.Cells.Item(11, 16).Value = 1600

$ws.Cells.Item(12, 4).Value = 44761
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 23000
$ws.Cells.Item(12, 12).Value = 25000
$ws.Cells.Item(12, 13).Value = 24000
$ws.Cells.Item(12, 16).Value = 1600

$ws.Cells.Item(13, 4).Value = 44781
$ws.Cells.Item(13, 10).Value = 70
$ws.Cells.Item(13, 11).Value = 24000
$ws.Cells.Item(13, 12).Value = 24000
$ws.Cells.Item(13, 13).Value = 24000
$ws.Cells.Item(13, 16).Value = 1600

$ws.Cells.Item(14, 4).Value = 44365
$ws.Cells.Item(14, 10).Value = 80
$ws.Cells.Item(14, 11).Value = 25000
$ws.Cells.Item(14, 12).Value = 25000
$ws.Cells.Item(14, 13).Value = 25000
$ws.Cells.Item(14, 16).Value = 1667

$ws.Cells.Item(15, 4).Value = 44771
$ws.Cells.Item(15, 10).Value = 90
$ws.Cells.Item(15, 11).Value = 25000
$ws.Cells.Item(15, 12).Value = 25000
$ws.Cells.Item(15, 13).Value = 25000
$ws.Cells.Item(15, 16).Value = 1667

$ws.Cells.Item(16, 4).Value = 44782
$ws.Cells.Item(16, 10).Value = 120
$ws.Cells.Item(16, 11).Value = 24000
$ws.Cells.Item(16, 12).Value = 24000
$ws.Cells.Item(16, 13).Value = 24000
$ws.Cells.Item(16, 16).Value = 1600
